$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.452.55"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.82%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.303.22"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.62%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.86%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.24"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.71%  "
$ws.Range("E7").Value = "  +1.42%  "
$ws.Range("E8").Value = "  +0.21%  "
$ws.Range("E9").Value = "  +0.08%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.09"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.82%  "
$ws.Range("E11").Value = "  -0.57%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.35"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.57%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.107"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.49%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.969"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.50%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.32"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.649.66"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.68%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.304.92"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.54%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.526.49"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.93%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000105"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.20%  "
$ws.Range("E21").Value = "  -1.99%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.56"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "276.63"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +6.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.13"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +19.68%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.29"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.72%  "
$ws.Range("E26").Value = "  -0.22%  "
$ws.Range("E27").Value = "  -1.68%  "
$ws.Range("E28").Value = "  +5.60%  "
$ws.Range("E29").Value = "  -0.24%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.80"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.11%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "164.92"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.16%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0876"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.53%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.90"
$ws.Range("D33").Style = "Normal"
$ws.Range("E34").Value = "  +5.24%  "
$ws.Range("E35").Value = "  -10.55%  "
$ws.Range("E36").Value = "  -2.25%  "
$ws.Range("E37").Value = "  +5.46%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.59"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.16%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.74"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.46%  "
$ws.Range("E40").Value = "  -1.32%  "
$ws.Range("E41").Value = "  +2.42%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "69.96"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.70%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "95.21"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.53%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.228"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.18%  "
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "82.81"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +10.78%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.06"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.91%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "113.02"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.94%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.89"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.27%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.590.72"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.48%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0999"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.35%  "
